$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 89; existing rows 89..141 shift down to 90..142.
$ws.Rows("89:89").Insert()

# Populate the newly inserted row 89 with the new data point.
$ws.Cells.Item(89, 1).Value = 4
$ws.Cells.Item(89, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(89, 3).Value = "Los Lagos"
$ws.Cells.Item(89, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(89, 5).Value = 10
$ws.Cells.Item(89, 6).Value = 100112021
$ws.Cells.Item(89, 7).Value = "Ají"
$ws.Cells.Item(89, 8).Value = "Inferno"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 25
$ws.Cells.Item(89, 11).Value = 45000
$ws.Cells.Item(89, 12).Value = 45000
$ws.Cells.Item(89, 13).Value = 45000
$ws.Cells.Item(89, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(89, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value = 3750
$ws.Cells.Item(89, 17).Value = 12
$ws.Cells.Item(89, 18).Value = "Hortaliza"

Write-Output "row inserted"
